$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A7: rename variable name to clarify it's a boolean
$ws.Range("A7").Value = "IsCursorOnPrimaryDisplayBool"

# D8: change "Dispose" wording to "Stop / Reset", keeping bold runs intact.
$d8Text = "If CursorTrackerTimer not started: ACTION: start the timer`nIf CursorTrackerTimer already started (else):`n - if CursorTrackerTimer > e.g. 2 mins, ACTION: dim and Stop / Reset the Timer`n - else (if CursorTrackerTimer < 2 mins), do not dim, (rerun the loop) "
$ws.Range("D8").Value = $d8Text
$ws.Range("D8").Characters(1, 35).Font.Bold = $false
$ws.Range("D8").Characters(36, 8).Font.Bold = $true
$ws.Range("D8").Characters(44, 102).Font.Bold = $false
$ws.Range("D8").Characters(146, 11).Font.Bold = $true
$ws.Range("D8").Characters(157, 99).Font.Bold = $false

# B7: append a follow-up note about checking after resetting the timer,
# keeping "ACTION: " bold and the rest regular.
$b7Text = "T: `nIf CursorTrackerTimer already started:`nACTION: reset the CursorTrackerTimer (cursor is back on the primary display). Then check:"
$ws.Range("B7").Value = $b7Text
$ws.Range("B7").Characters(1, 43).Font.Bold = $false
$ws.Range("B7").Characters(44, 8).Font.Bold = $true
$ws.Range("B7").Characters(52, 81).Font.Bold = $false

# Update the active selection to match the target view
$ws.Range("D16").Select()
